# Update the Config sheet: row 2 (first data row) now targets the
# LeaveDeduction module, run as user "bob" instead of "automation3".
$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Config")

$ws.Range("B2").Value = "mod:LeaveDeduction"
$ws.Range("D2").Value = "bob"

# Move the active selection to D2, matching the cell that was edited.
$ws.Range("D2").Select()
